# Weekly ranking update [2025-12-17]
# Adds a new worksheet "magapoke_2025-12-17" at the end of the workbook,
# populated with the rank/title table for that week, formatted to match
# the header style used on the other weekly sheets.

$wb = $excel.ActiveWorkbook

# Create the new sheet and name it before positioning it, so lookups by
# name are unambiguous for the rest of the script.
$ws = $wb.Worksheets.Add()
$ws.Name = "magapoke_2025-12-17"

# Header row
$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "title"

# Data rows: rank (number) + title (text)
$titles = @(
    "アイドラトリィ",
    "せいぶつ部の田辺くん",
    "黒月のイェルクナハト",
    "スルガメテオ",
    "ドリーム☆ジャンボ☆ガール",
    "黄昏町プリズナーズ",
    "K-9~警視庁公安部公安第9課異能対策係~",
    "篝家の８兄弟",
    "ハードワーカー中田",
    "生きたがりの人狼",
    "ナキナギ",
    "永久のユウグレ",
    "ルックスＹを選んでしまいました 〜やり込んでいるゲームに転生したはずなのに、未実装のガチャで攻略をすることになった件〜",
    "普通の本はありません！",
    "平成転生",
    "ゼロとヒャク",
    "お母さん冒険者、ログインボーナスでスキル【主婦】に目覚めました。週一貰えるチラシで冒険者生活頑張ります！",
    "鳴るさんだぁ",
    "夜鐘のキト",
    "春くらり",
    "その青春",
    "卒業アルバムの彼女たち",
    "屋根の下のアルテミス",
    "歪みの虜",
    "MYS",
    "ハナバス　苔石花江のバスケ論",
    "ともだちづくり",
    "異世界グルメで成り上がり無双～山に追放されたので、のんびりキャンプを楽しんでいたらいつの間にか強くなっていて、王侯貴族や実力者たちが俺を放っておいてくれません。一方、俺を追放した貴族たちは破滅が始まる～",
    "君が監督！",
    "白銀のキュイジーヌ～明治外交官の料理人～",
    "明智ナンバーワン",
    "皇女転生　～伝説の大魔導士（♂）、姫騎士となりて伝説の令嬢騎士団を作り無双する～",
    "ハプスブルク家の華麗なる受難",
    "JK Biker",
    "追放されなかった男　～二度目の人生は土下座から始まりました～",
    "限界集落を脱村した錬金術士、都会で`"最強`"なのがバレまくる。～老害どもにはいい加減愛想が尽きました～",
    "イエティ、とある日々",
    "眠れる森のレガ",
    "東京デスレース",
    "人生逆転ダンジョン",
    "じゅーくぼっくす",
    "〈小市民〉 春期限定いちごタルト事件",
    "ナマイキ旭ちゃんをわからせたい",
    "花子狩り"
)

for ($i = 0; $i -lt $titles.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $titles[$i]
}

# Header styling: bold font, thin box border, centered horizontally,
# top-aligned vertically -- matches the style used on the other sheets.
$header = $ws.Range("A1:B1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2

# Match the page margins used by the other weekly sheets.
$ws.PageSetup.LeftMargin = 0.75 * 72
$ws.PageSetup.RightMargin = 0.75 * 72
$ws.PageSetup.TopMargin = 1 * 72
$ws.PageSetup.BottomMargin = 1 * 72
$ws.PageSetup.HeaderMargin = 0.5 * 72
$ws.PageSetup.FooterMargin = 0.5 * 72

# Move the new sheet to the end of the tab order, after the most recent
# existing weekly sheet. Re-fetch by name afterwards since moving a sheet
# can invalidate earlier object references.
$ws.Move($null, $wb.Worksheets.Item("magapoke_2025-12-10"))

$ws = $wb.Worksheets.Item("magapoke_2025-12-17")
$ws.Range("A1").Select()
